$wb = $excel.ActiveWorkbook

# Incidental view-state leftover on the "Sheet1" tab (physically sheet2.xml)
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("B34").Select()

# Main edit happens on the "Sheet2" tab (physically sheet1.xml, the active tab)
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Activate()

$ws2.Range("A8").Value = "Visual"
$ws2.Range("B8").Value = "Visual Basic"
$ws2.Range("C8").Value = "V: select at line level`nv: select at char level`nCtrl + v : select vertically at char level"
$ws2.Range("C8").WrapText = $true
$ws2.Rows.Item(8).RowHeight = 45

$ws2.Range("C8").Select()
